$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.765.75'
$ws.Range('D3').Value = '1.537.95'
$ws.Range('E3').Value = '  -1.76%  '
$ws.Range('D5').Value = '205.66'
$ws.Range('E6').Value = '  -0.76%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D9').Value = '21.22'
$ws.Range('E9').Value = '  -2.82%  '
$ws.Range('E10').Value = '  -0.56%  '
$ws.Range('E11').Value = '  -1.10%  '
$ws.Range('D12').Value = '1.756.93'
$ws.Range('E12').Value = '  -1.76%  '
$ws.Range('D13').Value = '1.546.05'
$ws.Range('E13').Value = '  -1.21%  '
$ws.Range('E14').Value = '  -1.21%  '
$ws.Range('E15').Value = '  -1.05%  '
$ws.Range('D16').Value = '26.765.87'
$ws.Range('E16').Value = '  -0.08%  '
$ws.Range('D17').Value = '61.02'
$ws.Range('E17').Value = '  -0.47%  '
$ws.Range('D18').Value = '213.31'
$ws.Range('E18').Value = '  -0.79%  '
$ws.Range('E19').Value = '  -1.64%  '
$ws.Range('E20').Value = '  +0.87%  '
$ws.Range('D22').Value = '4.00'
$ws.Range('E22').Value = '  -2.04%  '
$ws.Range('E23').Value = '  -1.63%  '
$ws.Range('E24').Value = '  -3.38%  '
$ws.Range('D25').Value = '151.67'
$ws.Range('E25').Value = '  -0.58%  '
$ws.Range('D26').Value = '6.57'
$ws.Range('E26').Value = '  -2.22%  '
$ws.Range('E27').Value = '  -1.02%  '
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('E29').Value = '  -0.85%  '
$ws.Range('E30').Value = '  -1.01%  '
$ws.Range('E31').Value = '  -1.37%  '
$ws.Range('E32').Value = '  +2.12%  '
$ws.Range('D33').Value = '1.363.38'
$ws.Range('E34').Value = '  +0.23%  '
$ws.Range('E35').Value = '  -1.77%  '
$ws.Range('D36').Value = '0.957'
$ws.Range('E36').Value = '  +3.11%  '
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('E38').Value = '  +1.12%  '
$ws.Range('E39').Value = '  -1.24%  '
$ws.Range('E41').Value = '  -1.77%  '
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('D43').Value = '2.20'
$ws.Range('E43').Value = '  +0.71%  '
$ws.Range('E44').Value = '  -0.55%  '
$ws.Range('E45').Value = '  -3.01%  '
$ws.Range('D46').Value = '1.671.37'
$ws.Range('E46').Value = '  -1.73%  '
$ws.Range('D47').Value = '84.09'
$ws.Range('E47').Value = '  -1.87%  '
$ws.Range('E48').Value = '  +3.60%  '
$ws.Range('D49').Value = '0.0₇0966'
$ws.Range('E49').Value = '  -1.80%  '
$ws.Range('D50').Value = '0.0941'
$ws.Range('E50').Value = '  -0.81%  '
$ws.Range('E51').Value = '  +0.13%  '
